$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accident records (rows 104-113)
# Row 104
$ws.Range("A104").Value = 102
$ws.Range("B104").Value = 2
$ws.Range("C104").Value = "angled t-bone"
$ws.Range("D104").Value = "Yes"
$ws.Range("E104").Value = "Yes"
$ws.Range("F104").Value = "No"
$ws.Range("G104").Value = "Yes"
$ws.Range("H104").Value = "bus, sedan"
$ws.Range("I104").Value = "No"
$ws.Range("J104").Value = "Cloudy"
$ws.Range("K104").Value = "No"
$ws.Range("L104").Value = "super busy intersection"

# Row 105
$ws.Range("A105").Value = 103
$ws.Range("B105").Value = 2
$ws.Range("C105").Value = "rear-end"
$ws.Range("D105").Value = "No"
$ws.Range("E105").Value = "No"
$ws.Range("F105").Value = "No"
$ws.Range("G105").Value = "No"
$ws.Range("H105").Value = "sedan x2"
$ws.Range("I105").Value = "No"
$ws.Range("J105").Value = "Day, clear"
$ws.Range("K105").Value = "No"
$ws.Range("L105").Value = "super busy intersection"

# Row 106
$ws.Range("A106").Value = 104
$ws.Range("B106").Value = 2
$ws.Range("C106").Value = "angled t-bone"
$ws.Range("D106").Value = "Yes"
$ws.Range("E106").Value = "Yes"
$ws.Range("F106").Value = "No"
$ws.Range("G106").Value = "Yes"
$ws.Range("H106").Value = "suv, sedan"
$ws.Range("I106").Value = "No"
$ws.Range("J106").Value = "Night"
$ws.Range("K106").Value = "Yes"
$ws.Range("L106").Value = "extremely difficult to see"

# Row 107
$ws.Range("A107").Value = 105
$ws.Range("B107").Value = 2
$ws.Range("C107").Value = "t-bone"
$ws.Range("D107").Value = "Yes"
$ws.Range("E107").Value = "Yes"
$ws.Range("F107").Value = "No"
$ws.Range("G107").Value = "Yes"
$ws.Range("H107").Value = "bike, suv"
$ws.Range("I107").Value = "No"
$ws.Range("J107").Value = "Day"
$ws.Range("K107").Value = "No"
$ws.Range("L107").Value = "bike runs into suv, only biker seems injured"

# Row 108
$ws.Range("A108").Value = 106
$ws.Range("B108").Value = 2
$ws.Range("C108").Value = "t-bone"
$ws.Range("D108").Value = "Yes"
$ws.Range("E108").Value = "Yes"
$ws.Range("F108").Value = "No"
$ws.Range("G108").Value = "Yes"
$ws.Range("H108").Value = "bike, truck"
$ws.Range("I108").Value = "No"
$ws.Range("J108").Value = "Day"
$ws.Range("K108").Value = "No"
$ws.Range("L108").Value = "Truck runs into bike"

# Row 109
$ws.Range("A109").Value = 107
$ws.Range("B109").Value = 3
$ws.Range("C109").Value = "rear-end"
$ws.Range("D109").Value = "Yes"
$ws.Range("E109").Value = "Yes"
$ws.Range("F109").Value = "No"
$ws.Range("G109").Value = "Yes"
$ws.Range("H109").Value = "sedan x2, bus"
$ws.Range("I109").Value = "No"
$ws.Range("J109").Value = "Day"
$ws.Range("K109").Value = "No"
$ws.Range("L109").Value = "Bus runs into sedan which runs into other sedan"

# Row 110
$ws.Range("A110").Value = 108
$ws.Range("B110").Value = 1
$ws.Range("C110").Value = "loss of control"
$ws.Range("D110").Value = "Yes"
$ws.Range("E110").Value = "Yes"
$ws.Range("F110").Value = "No"
$ws.Range("G110").Value = "Yes"
$ws.Range("H110").Value = "sedan"
$ws.Range("I110").Value = "No"
$ws.Range("J110").Value = "Day"
$ws.Range("K110").Value = "No"
$ws.Range("L110").Value = "car runs into median in city"

# Row 111
$ws.Range("A111").Value = 109
$ws.Range("B111").Value = 2
$ws.Range("C111").Value = "t-bone"
$ws.Range("D111").Value = "Yes"
$ws.Range("E111").Value = "Yes"
$ws.Range("F111").Value = "No"
$ws.Range("G111").Value = "Yes"
$ws.Range("H111").Value = "sedan, bus"
$ws.Range("I111").Value = "No"
$ws.Range("J111").Value = "Night"
$ws.Range("K111").Value = "No"
$ws.Range("L111").Value = "sedan runs into bus at non-busy intersection"

# Row 112
$ws.Range("A112").Value = 110
$ws.Range("B112").Value = 1
$ws.Range("C112").Value = "topple"
$ws.Range("D112").Value = "Yes"
$ws.Range("E112").Value = "Yes"
$ws.Range("F112").Value = "No"
$ws.Range("G112").Value = "Yes"
$ws.Range("H112").Value = "bike"
$ws.Range("I112").Value = "No"
$ws.Range("J112").Value = "Day"
$ws.Range("K112").Value = "No"
$ws.Range("L112").Value = "Bike topples over, injury seems minor, unknown if any other vehicle was involved"

# Row 113
$ws.Range("A113").Value = 111
$ws.Range("B113").Value = 2
$ws.Range("C113").Value = "t-bone"
$ws.Range("D113").Value = "Yes"
$ws.Range("E113").Value = "Yes"
$ws.Range("F113").Value = "No"
$ws.Range("G113").Value = "Yes"
$ws.Range("H113").Value = "van x2"
$ws.Range("I113").Value = "No"
$ws.Range("J113").Value = "Night"
$ws.Range("K113").Value = "Yes"
$ws.Range("L113").Value = "The glare makes it difficult to tell what exactly happened"

# Update view state: scroll position and final selection
[void]$excel.Goto($ws.Range("A77"), $true)
[void]$ws.Range("B114").Select()
